$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) on cells whose new numeric-looking value
# would otherwise be auto-converted/reformatted by Excel (e.g. losing trailing
# zeros or switching to scientific notation), so the literal string is preserved
# exactly as in the source data feed.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume figures.
$ws.Range("D2").Value = "28.414.39"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "314.34"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.5453"
$ws.Range("E7").Value = "  +4.44%  "
$ws.Range("D8").Value = "0.3823"
$ws.Range("E8").Value = "  +4.35%  "
$ws.Range("D9").Value = "0.07585"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "42.51"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "1.120"
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "6.189"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "7.391"
$ws.Range("E15").Value = "  +6.60%  "
$ws.Range("D16").Value = "1.795.38"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").Value = "91.38"
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").Value = "0.00001070"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "0.06460"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").Value = "5.957"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").Value = "28.410.55"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("D24").Value = "11.40"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "2.120"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "160.18"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").Value = "20.72"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "2.393"
$ws.Range("E28").Value = "  +3.45%  "
$ws.Range("D29").Value = "2.000.63"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "123.12"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").Value = "1.125"
$ws.Range("E31").Value = "  +6.84%  "
$ws.Range("D32").Value = "0.1026"
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("D33").Value = "5.740"
$ws.Range("E33").Value = "  +3.41%  "
$ws.Range("D34").Value = "3.682"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "0.2315"
$ws.Range("E35").Value = "  +15.00%  "
$ws.Range("D36").Value = "0.06546"
$ws.Range("E36").Value = "  +10.00%  "
$ws.Range("D37").Value = "0.02319"
$ws.Range("D38").Value = "5.189"
$ws.Range("E38").Value = "  +7.52%  "
$ws.Range("D39").Value = "8.754"
$ws.Range("E39").Value = "  +8.84%  "
$ws.Range("D40").Value = "11.62"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("E41").Value = "  +4.38%  "
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "1.390"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").Value = "13.56"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("D46").Value = "0.5953"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").Value = "3.676"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "126.25"
$ws.Range("E48").Value = "  +4.31%  "
$ws.Range("D49").Value = "1.993"
$ws.Range("E49").Value = "  +6.45%  "
$ws.Range("D50").Value = "1.151"
$ws.Range("E50").Value = "  +3.36%  "
